# Updated cryptos list on Fri Jan 12 15:39:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.783.10"
$ws.Range("E2").Value = "'  -6.06%  "

$ws.Range("D3").Value = "'2.659.42"
$ws.Range("E3").Value = "'  +0.82%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'306.27"
$ws.Range("E5").Value = "'  -1.83%  "

$ws.Range("D6").Value = "'96.86"
$ws.Range("E6").Value = "'  -6.28%  "

$ws.Range("D7").Value = "'0.589"
$ws.Range("E7").Value = "'  -3.30%  "

$ws.Range("E8").Value = "'  +0.12%  "

$ws.Range("D9").Value = "'0.570"
$ws.Range("E9").Value = "'  -4.42%  "

$ws.Range("D10").Value = "'37.67"
$ws.Range("E10").Value = "'  -6.73%  "

$ws.Range("D11").Value = "'0.0835"
$ws.Range("E11").Value = "'  -3.23%  "

$ws.Range("D12").Value = "'7.96"
$ws.Range("E12").Value = "'  -4.51%  "

$ws.Range("D13").Value = "'3.087.84"
$ws.Range("E13").Value = "'  +1.44%  "

$ws.Range("E14").Value = "'  -0.06%  "

$ws.Range("D15").Value = "'2.678.33"
$ws.Range("E15").Value = "'  +1.28%  "

$ws.Range("D16").Value = "'0.915"
$ws.Range("E16").Value = "'  -2.19%  "

$ws.Range("D17").Value = "'14.93"
$ws.Range("E17").Value = "'  -2.35%  "

$ws.Range("D18").Value = "'44.865.77"
$ws.Range("E18").Value = "'  -6.12%  "

$ws.Range("D19").Value = "'6.80"
$ws.Range("E19").Value = "'  -0.37%  "

$ws.Range("D20").Value = "'0.0₃0999"
$ws.Range("E20").Value = "'  -3.29%  "

$ws.Range("D21").Value = "'12.54"
$ws.Range("E21").Value = "'  -5.24%  "

$ws.Range("D22").Value = "'74.52"
$ws.Range("E22").Value = "'  +1.59%  "

$ws.Range("D23").Value = "'276.09"
$ws.Range("E23").Value = "'  -1.61%  "

$ws.Range("E24").Value = "'  +3.17%  "

$ws.Range("D25").Value = "'3.00"
$ws.Range("E25").Value = "'  -2.60%  "

$ws.Range("D26").Value = "'30.76"
$ws.Range("E26").Value = "'  +0.02%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "'  +0.03%  "

$ws.Range("D28").Value = "'10.42"
$ws.Range("E28").Value = "'  -2.68%  "

$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "'  -3.91%  "

$ws.Range("D30").Value = "'37.46"
$ws.Range("E30").Value = "'  -6.05%  "

$ws.Range("B31").Value = "'Filecoin"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'6.13"
$ws.Range("E31").Value = "'  -1.67%  "

$ws.Range("B32").Value = "'LidoDAOToken"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D32").Value = "'3.74"
$ws.Range("E32").Value = "'  +0.82%  "

$ws.Range("D33").Value = "'2.31"
$ws.Range("E33").Value = "'  +3.09%  "

$ws.Range("B34").Value = "'WEMIXToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.81"
$ws.Range("E34").Value = "'  -2.23%  "

$ws.Range("B35").Value = "'Monero"
$ws.Range("C35").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'152.84"
$ws.Range("E35").Value = "'  -0.23%  "

$ws.Range("D36").Value = "'0.0828"
$ws.Range("E36").Value = "'  -3.66%  "

$ws.Range("D37").Value = "'0.119"
$ws.Range("E37").Value = "'  -7.61%  "

$ws.Range("B38").Value = "'Stellar"
$ws.Range("C38").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.123"
$ws.Range("E38").Value = "'  -1.35%  "

$ws.Range("B39").Value = "'EnergySwap"
$ws.Range("C39").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'24.66"
$ws.Range("E39").Value = "'  +8.77%  "

$ws.Range("D40").Value = "'15.86"
$ws.Range("E40").Value = "'  -1.35%  "

$ws.Range("D41").Value = "'3.56"
$ws.Range("E41").Value = "'  -4.44%  "

$ws.Range("D42").Value = "'0.0320"
$ws.Range("E42").Value = "'  -4.71%  "

$ws.Range("D43").Value = "'2.135.87"
$ws.Range("E43").Value = "'  -1.78%  "

$ws.Range("D44").Value = "'3.91"
$ws.Range("E44").Value = "'  -8.35%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  +0.08%  "

$ws.Range("D46").Value = "'91.63"
$ws.Range("E46").Value = "'  -6.27%  "

$ws.Range("D47").Value = "'9.39"
$ws.Range("E47").Value = "'  -4.67%  "

$ws.Range("D48").Value = "'2.935.76"
$ws.Range("E48").Value = "'  +1.38%  "

$ws.Range("D49").Value = "'109.64"
$ws.Range("E49").Value = "'  -0.50%  "

$ws.Range("D50").Value = "'1.60"
$ws.Range("E50").Value = "'  -2.71%  "

$ws.Range("D51").Value = "'0.196"
$ws.Range("E51").Value = "'  -3.93%  "

Write-Host "Updated cryptos list"